# Update ticker symbols in columns B, C, D, F for rows 2-52
# per the refreshed watchlist data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NSE:63MOONS"
$ws.Range("C2").Value = "NSE:ARTNIRMAN"
$ws.Range("D2").Value = ""
$ws.Range("B3").Value = "NSE:ABAN"
$ws.Range("C3").Value = "NSE:AVROIND"
$ws.Range("D3").Value = ""
$ws.Range("F3").Value = "NSE:BIOCON"
$ws.Range("B4").Value = "NSE:ADVENZYMES"
$ws.Range("C4").Value = "NSE:DIVGIITTS"
$ws.Range("D4").Value = ""
$ws.Range("F4").Value = "NSE:DIVISLAB"
$ws.Range("B5").Value = "NSE:ALPA"
$ws.Range("C5").Value = "NSE:INDRAMEDCO"
$ws.Range("F5").Value = "NSE:GAIL"
$ws.Range("B6").Value = "NSE:ALPHAGEO"
$ws.Range("C6").Value = "NSE:LTGILTBEES"
$ws.Range("F6").Value = "NSE:HEROMOTOCO"
$ws.Range("B7").Value = "NSE:APEX"
$ws.Range("C7").Value = "NSE:MICEL"
$ws.Range("F7").Value = "NSE:INFY"
$ws.Range("B8").Value = "NSE:BAJAJHIND"
$ws.Range("C8").Value = "NSE:NDGL"
$ws.Range("B9").Value = "NSE:BALPHARMA"
$ws.Range("C9").Value = "NSE:ORIENTCEM"
$ws.Range("B10").Value = "NSE:BIOCON"
$ws.Range("B11").Value = "NSE:BSHSL"
$ws.Range("B12").Value = "NSE:BUTTERFLY"
$ws.Range("B13").Value = "NSE:CAMLINFINE"
$ws.Range("B14").Value = "NSE:CARERATING"
$ws.Range("B15").Value = "NSE:CHENNPETRO"
$ws.Range("B16").Value = "NSE:CINEVISTA"
$ws.Range("B17").Value = "NSE:CIPLA"
$ws.Range("B18").Value = "NSE:DCMSHRIRAM"
$ws.Range("B19").Value = "NSE:DELTACORP"
$ws.Range("B20").Value = "NSE:DHANI"
$ws.Range("B21").Value = "NSE:DMCC"
$ws.Range("B22").Value = "NSE:DWARKESH"
$ws.Range("B23").Value = "NSE:ELECTCAST"
$ws.Range("B24").Value = "NSE:EMIL"
$ws.Range("B25").Value = "NSE:ENERGYDEV"
$ws.Range("B26").Value = "NSE:FCL"
$ws.Range("B27").Value = "NSE:FCSSOFT"
$ws.Range("B28").Value = "NSE:GAIL"
$ws.Range("B29").Value = "NSE:GFLLIMITED"
$ws.Range("B30").Value = "NSE:HAVISHA"
$ws.Range("B31").Value = "NSE:HCLTECH"
$ws.Range("B32").Value = "NSE:HGS"
$ws.Range("B33").Value = "NSE:HYBRIDFIN"
$ws.Range("B34").Value = "NSE:INDOWIND"
$ws.Range("B35").Value = "NSE:IOLCP"
$ws.Range("B36").Value = "NSE:JAYBARMARU"
$ws.Range("B37").Value = "NSE:KNRCON"
$ws.Range("B38").Value = "NSE:KOPRAN"
$ws.Range("B39").Value = "NSE:LALPATHLAB"
$ws.Range("B40").Value = "NSE:LTTS"
$ws.Range("B41").Value = "NSE:LXCHEM"
$ws.Range("B42").Value = "NSE:MAHSEAMLES"
$ws.Range("B43").Value = "NSE:MANGALAM"
$ws.Range("B44").Value = "NSE:MAWANASUG"
$ws.Range("B45").Value = "NSE:MURUDCERA"
$ws.Range("B46").Value = "NSE:NATCOPHARM"
$ws.Range("B47").Value = "NSE:NATHBIOGEN"
$ws.Range("B48").Value = "NSE:NINSYS"
$ws.Range("B49").Value = "NSE:NSLNISP"
$ws.Range("B50").Value = "NSE:ONMOBILE"
$ws.Range("B51").Value = "NSE:PARACABLES"
$ws.Range("B52").Value = "NSE:PATELENG"

# Append 10 new watchlist rows (53-62), copying the formatting
# (bold, centered, bordered) already used for the row-index column.
$ws.Range("A52").Copy()
$ws.Range("A53:A62").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A53").Value = 51
$ws.Range("B53").Value = "NSE:PCJEWELLER"
$ws.Range("A54").Value = 52
$ws.Range("B54").Value = "NSE:PFS"
$ws.Range("A55").Value = 53
$ws.Range("B55").Value = "NSE:RAIN"
$ws.Range("A56").Value = 54
$ws.Range("B56").Value = "NSE:RANASUG"
$ws.Range("A57").Value = 55
$ws.Range("B57").Value = "NSE:RENUKA"
$ws.Range("A58").Value = 56
$ws.Range("B58").Value = "NSE:RGL"
$ws.Range("A59").Value = 57
$ws.Range("B59").Value = "NSE:RTNINDIA"
$ws.Range("A60").Value = 58
$ws.Range("B60").Value = "NSE:RTNPOWER"
$ws.Range("A61").Value = 59
$ws.Range("B61").Value = "NSE:RUPA"
$ws.Range("A62").Value = 60
$ws.Range("B62").Value = "NSE:SAKSOFT"
